$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3494.8333
$ws.Range("I7").Value = 3942.25
$ws.Range("J7").Value = 2600
$ws.Range("K7").Value = 3942.25
$ws.Range("L7").Value = 2600
$ws.Range("M7").Value = -3830.25
$ws.Range("N7").Value = -2824

$ws.Range("H14").Value = 3494.8333
$ws.Range("I14").Value = 3942.25
$ws.Range("J14").Value = 2600
$ws.Range("K14").Value = 3942.25
$ws.Range("L14").Value = 2600
$ws.Range("M14").Value = -3751.25
$ws.Range("N14").Value = -2982

$ws.Range("H70").Value = 6783.1665
$ws.Range("J70").Value = 9250
$ws.Range("L70").Value = 27750
$ws.Range("N70").Value = -28290

$ws.Range("H73").Value = 6783.1665
$ws.Range("J73").Value = 9250
$ws.Range("L73").Value = 27750
$ws.Range("N73").Value = -29622

$ws.Range("H96").Value = 677.35297
$ws.Range("I96").Value = 874.1
$ws.Range("J96").Value = 396.2857
$ws.Range("K96").Value = 2622.3
$ws.Range("L96").Value = 1188.8571
$ws.Range("M96").Value = -1249.3
$ws.Range("N96").Value = -3934.8571

$ws.Range("H98").Value = 1298.1666
$ws.Range("I98").Value = 1298.1666
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1298.1666
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 199.8334
$ws.Range("N98").ClearContents()

$ws.Range("H112").Value = 202969.4
$ws.Range("I112").Value = 201219
$ws.Range("J112").Value = 204719.8
$ws.Range("K112").Value = 603657
$ws.Range("L112").Value = 614159.3999999999
$ws.Range("M112").Value = -602549
$ws.Range("N112").Value = -616375.3999999999

$ws.Range("H115").Value = 1052.25
$ws.Range("I115").Value = 569.6667
$ws.Range("K115").Value = 1709.0001
$ws.Range("M115").Value = -142.0001

$ws.Range("H122").Value = 1298.1666
$ws.Range("I122").Value = 1298.1666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3894.4998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1444.4998
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2869.0303
$ws.Range("I132").Value = 2914.9375
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 8744.8125
$ws.Range("L132").Value = 4200
$ws.Range("M132").Value = -6214.8125
$ws.Range("N132").Value = -9260

$ws.Range("H137").Value = 848.06665
$ws.Range("I137").Value = 697.7692
$ws.Range("K137").Value = 2093.3076
$ws.Range("M137").Value = 456.6923999999999

$ws.Range("H138").Value = 3251.5952
$ws.Range("J138").Value = 4339.04
$ws.Range("L138").Value = 13017.12
$ws.Range("N138").Value = -23297.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1226233.1
$ws.Range("I2").Value = 1337345.1
$ws.Range("K2").Value = 1337345.1
$ws.Range("M2").Value = -1337232.1

$ws.Range("H61").Value = 83338776
$ws.Range("I61").Value = 142859900
$ws.Range("K61").Value = 142859900
$ws.Range("M61").Value = -142859688

$ws.Range("H74").Value = 23811216
$ws.Range("I74").Value = 31251600
$ws.Range("J74").Value = 1992.9
$ws.Range("K74").Value = 31251600
$ws.Range("L74").Value = 1992.9
$ws.Range("M74").Value = -31250726
$ws.Range("N74").Value = -3740.9

$ws.Range("H77").Value = 23811216
$ws.Range("I77").Value = 31251600
$ws.Range("J77").Value = 1992.9
$ws.Range("K77").Value = 156258000
$ws.Range("L77").Value = 9964.5
$ws.Range("M77").Value = -156253632
$ws.Range("N77").Value = -18700.5

$ws.Range("H116").Value = 1226233.1
$ws.Range("I116").Value = 1337345.1
$ws.Range("K116").Value = 1337345.1
$ws.Range("M116").Value = -1335051.1

$ws.Range("H136").Value = 83338776
$ws.Range("I136").Value = 142859900
$ws.Range("K136").Value = 428579700
$ws.Range("M136").Value = -428577150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1226233.1
$ws.Range("I3").Value = 1337345.1
$ws.Range("K3").Value = 1337345.1
$ws.Range("M3").Value = -1337231.1

$ws.Range("H20").Value = 2253
$ws.Range("I20").Value = 2188.1333
$ws.Range("K20").Value = 2188.1333
$ws.Range("M20").Value = -1941.1333

$ws.Range("H108").Value = 71165.664
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 71165.664
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 71165.664
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -78845.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1466.6666
$ws.Range("I15").Value = 1400
$ws.Range("K15").Value = 1400
$ws.Range("M15").Value = -1230

$ws.Range("H132").Value = 41669684
$ws.Range("I132").Value = 41669684
$ws.Range("K132").Value = 125009052
$ws.Range("M132").Value = -125006522

$ws.Range("H134").Value = 14707611
$ws.Range("I134").Value = 14707611
$ws.Range("K134").Value = 44122833
$ws.Range("M134").Value = -44120298

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 200
$ws.Range("K9").Value = 600
$ws.Range("M9").Value = -376

$ws.Range("H34").Value = 929.35297
$ws.Range("J34").Value = 1000.0333
$ws.Range("L34").Value = 3000.0999
$ws.Range("N34").Value = -3168.0999

$ws.Range("H94").Value = 12257.077
$ws.Range("I94").Value = 1498.5
$ws.Range("J94").Value = 21478.715
$ws.Range("K94").Value = 4495.5
$ws.Range("L94").Value = 64436.145
$ws.Range("M94").Value = -3819.5
$ws.Range("N94").Value = -65788.145

$ws.Range("H107").Value = 1274.75
$ws.Range("I107").Value = 334.7
$ws.Range("K107").Value = 1004.1
$ws.Range("M107").Value = 915.9000000000001

$ws.Range("H113").Value = 67308.92999999999
$ws.Range("I113").Value = 166949.17
$ws.Range("K113").Value = 500847.51
$ws.Range("M113").Value = -498677.51

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 15564.375
$ws.Range("J12").Value = 100000
$ws.Range("L12").Value = 100000
$ws.Range("N12").Value = -100280

$ws.Range("H70").Value = 5375.6
$ws.Range("I70").Value = 4969.625
$ws.Range("K70").Value = 4969.625
$ws.Range("M70").Value = -4699.625

$ws.Range("H73").Value = 5375.6
$ws.Range("I73").Value = 4969.625
$ws.Range("K73").Value = 4969.625
$ws.Range("M73").Value = -4033.625

$ws.Range("H102").Value = 3298
$ws.Range("I102").Value = 3298
$ws.Range("K102").Value = 3298
$ws.Range("M102").Value = -1676

$ws.Range("H122").Value = 2809.3333
$ws.Range("I122").Value = 1733.8695
$ws.Range("J122").Value = 5282.9
$ws.Range("K122").Value = 5201.6085
$ws.Range("L122").Value = 15848.7
$ws.Range("M122").Value = -2751.6085
$ws.Range("N122").Value = -20748.7

$ws.Range("H132").Value = 5954549.5
$ws.Range("I132").Value = 6946363.5
$ws.Range("K132").Value = 20839090.5
$ws.Range("M132").Value = -20836560.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 576
$ws.Range("I55").Value = 574.3570999999999
$ws.Range("J55").Value = 581.75
$ws.Range("K55").Value = 574.3570999999999
$ws.Range("L55").Value = 581.75
$ws.Range("M55").Value = -401.3570999999999
$ws.Range("N55").Value = -927.75

$ws.Range("H68").Value = 5901.4
$ws.Range("I68").Value = 1821.25
$ws.Range("J68").Value = 22222
$ws.Range("K68").Value = 1821.25
$ws.Range("L68").Value = 22222
$ws.Range("M68").Value = -1072.25
$ws.Range("N68").Value = -23720

$ws.Range("H71").Value = 5901.4
$ws.Range("I71").Value = 1821.25
$ws.Range("J71").Value = 22222
$ws.Range("K71").Value = 9106.25
$ws.Range("L71").Value = 111110
$ws.Range("M71").Value = -5362.25
$ws.Range("N71").Value = -118598

$ws.Range("H96").Value = 32499
$ws.Range("J96").Value = 32499
$ws.Range("L96").Value = 32499
$ws.Range("N96").Value = -37991

$ws.Range("H122").Value = 3237
$ws.Range("I122").Value = 3237
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9711
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7261
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H132").Value = 31250886
$ws.Range("I132").Value = 38462356
$ws.Range("K132").Value = 115387068
$ws.Range("M132").Value = -115384538
